$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting existing rows 56..152 down to 57..153
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly data entry
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44915
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = 100114002
$ws.Range("G56").Value = "Camote"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 20
$ws.Range("K56").Value = 24000
$ws.Range("L56").Value = 24000
$ws.Range("M56").Value = 24000
$ws.Range("N56").Value = "$/malla 20 kilos"
$ws.Range("O56").Value = "Perú"
$ws.Range("P56").Value = 1200
$ws.Range("Q56").Value = 20
$ws.Range("R56").Value = "Hortaliza"
